$wb = $excel.ActiveWorkbook
$details = $wb.Worksheets.Item("Details")

# ---- Fill in the candidate's details in column C of the "Details" sheet ----
$details.Cells.Item(5, 3).Value  = "JONES"                # Last Name
$details.Cells.Item(6, 3).Value  = "David"                # First Name

$dob = Get-Date -Year 1947 -Month 1 -Day 8 -Hour 0 -Minute 0 -Second 0
$details.Cells.Item(8, 3).Value  = $dob                   # Date of Birth

$details.Cells.Item(9, 3).Value  = "Chateau de Signal"    # Address line 1
$details.Cells.Item(11, 3).Value = "Blonay"               # City
$details.Cells.Item(12, 3).Value = "VD"                   # Canton
$details.Cells.Item(13, 3).Value = "Switzerland"          # Country
$details.Cells.Item(14, 3).Value = "Married"              # Marital status
$details.Cells.Item(15, 3).Value = 3                      # Count of children
$details.Cells.Item(16, 3).Value = "No"                   # PhD?
$details.Cells.Item(18, 3).Value = "MA"                   # Highest academic credential (if no PhD)

$diploma = Get-Date -Year 1963 -Month 7 -Day 5 -Hour 0 -Minute 0 -Second 0
$details.Cells.Item(19, 3).Value = $diploma                # Date of Diploma

$details.Cells.Item(20, 3).Value = "Music"                # Profession/Academic Focus
$details.Cells.Item(21, 3).Value = "U.K."                 # Nationality
$details.Cells.Item(23, 3).Value = "YES"                  # Work permit?
$details.Cells.Item(24, 3).Value = "C"                    # Type of work permit

$startDate = Get-Date -Year 2016 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$details.Cells.Item(25, 3).Value = $startDate               # Start date availability

$details.Cells.Item(26, 3).Value = "Musician"              # Role
$details.Cells.Item(27, 3).Value = "R"                     # Team

$pct = $details.Cells.Item(28, 3)
$pct.Value = 1                                              # % worked
$pct.NumberFormat = "0%"

$details.Cells.Item(29, 3).Value = "Lausanne"               # Place of work

# ---- Update the view state: Details becomes the active tab, selection moves to B29 ----
$cdc = $wb.Worksheets.Item("CdC")
$cdc.Range("B11").Select()

$prop = $wb.Worksheets.Item("Prop d'engagement")
$prop.Range("D25").Select()

$details.Activate()
$details.Range("B29").Select()
